$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F ("想去人数") that needs updating.
$updates = @{
    3  = 93
    4  = 1541
    5  = 588
    6  = 1083
    7  = 11217
    8  = 6
    9  = 86
    10 = 95
    11 = 333
    12 = 1076
    14 = 12280
    15 = 12886
    17 = 132
    19 = 31
    22 = 66
}

# Both "展览" (sheet 1) and "全部类型" (sheet 4) contain the same data table
# and need the same updates applied to column F.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
